# Update Turk input file: drop the old "video_url"/Hit*_video_url_data
# column (B) entirely, shifting the former image2_url data (column C)
# left into column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column B (video_url + Hit*_video_url_data values). Column C
# (image2_url + its video2frame... urls) shifts left and becomes the
# new column B.
$ws.Columns("B").Delete()

# Select column B like the saved file does.
$ws.Range("B1:B1048576").Select()
